$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1 — new columns for 09-04-2025 and 10-04-2025
$ws.Range("F1").Value = "09-04-2025 Status"
$ws.Range("G1").Value = "09-04-2025 Time"
$ws.Range("H1").Value = "10-04-2025 Status"
$ws.Range("I1").Value = "10-04-2025 Time"

# Copy the header style from D1 (existing "Status/Time" header) onto the new headers
$ws.Range("D1:E1").Copy()
$ws.Range("F1:I1").PasteSpecial(-4122)

# Fill data rows 2-42 with Status = "A" and Time = "00:00:00"
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 6).Value = "A"
    $ws.Cells.Item($r, 7).Value = "00:00:00"
    $ws.Cells.Item($r, 8).Value = "A"
    $ws.Cells.Item($r, 9).Value = "00:00:00"
}
